# "If Stock dispo alors réduction ok + expédié status ok"
# Renumber the order numbers in column A (Numéro de commande) from the
# placeholder 1001-1005 values down to a simple sequential 1-5, and move
# the active selection/cursor the way the author left it (H7) on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# Restore the view/selection state saved with the workbook.
[void]$ws.Range("H7").Select()
